$d = $word.ActiveDocument

# 1) Expand the "zipped archive" parenthetical in the submission instructions.
$d.Content.Find.Execute(
    "(zipped archive) on Canvas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(consisting of a PDF report and a zipped archive containing additional files, e.g., test and model code) on Canvas.",
    2
)

# 2) Report-template note gains a sentence about adapting it freely.
$d.Content.Find.Execute(
    "A template for your report can be found here.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A template for your report can be found here. You may modify this template for your own purposes, if needed.",
    2
)

# 3) CoffeeMaker question now references Assignment 3 instead of Assignment 2.
$d.Content.Find.Execute(
    "CoffeeMaker example from Assignment 2.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CoffeeMaker example from Assignment 3.",
    2
)

# 4) Bold "Assignment 2" reference inside the table also becomes "Assignment 3".
$d.Content.Find.Execute(
    "Assignment 2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Assignment 3",
    2
)

Write-Output "done"
